$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# "Contestants" sheet: Kathleen Reynolds (row 3) and Felicity Parker-Hill
# (row 4) swap places. Felicity now lands on row 3 (status becomes
# "available", and she carries an explicit-but-empty Rating value), and
# Kathleen now lands on row 4 (status stays "assigned", Rating cleared).
# ---------------------------------------------------------------------------
$contestants = $wb.Worksheets.Item("Contestants")

# Row 3 -> Felicity Parker-Hill
$contestants.Range("A3").Value = "28603f95-d5f6-47ab-88c4-0d79742a6b02"
$contestants.Range("B3").Value = "Felicity Parker-Hill"
$contestants.Range("C3").Value = 27
$contestants.Range("E3").Value = "felicity.parkerhill@endemolshine.com.au"
$contestants.Range("G3").Value = "Melbourne"
$contestants.Range("H3").Value = ""
$contestants.Range("I3").Value = "available"
$contestants.Range("J3").Value = "Peter Adamidis, Kathleen Reynolds"

# Row 4 -> Kathleen Reynolds
$contestants.Range("A4").Value = "d698b1de-6641-45c6-aa63-f577d2b634bb"
$contestants.Range("B4").Value = "Kathleen Reynolds"
$contestants.Range("C4").Value = 33
$contestants.Range("E4").Value = "kathleenmonicareynolds@gmail.com"
$contestants.Range("G4").Value = "Footscray"
$contestants.Range("H4").Value = ""
$contestants.Range("I4").Value = "assigned"
$contestants.Range("J4").Value = "Peter Adamidis, Felicity Parker-Hill"

# ---------------------------------------------------------------------------
# "Seat Assignments" sheet: Felicity's seat-assignment row is removed, and
# Kathleen's row is consolidated into row 2 with a new id / block / seat.
# ---------------------------------------------------------------------------
$seats = $wb.Worksheets.Item("Seat Assignments")

# Drop the old row 3 (Kathleen's original assignment record) entirely -
# row 2 (Felicity's) shifts up and gets overwritten below.
$seats.Rows.Item(3).Delete()

$seats.Range("A2").Value = "6da0092e-89ab-452b-80c7-216b45398ce1"
$seats.Range("C2").Value = "d698b1de-6641-45c6-aa63-f577d2b634bb"
$seats.Range("D2").Value = 2
$seats.Range("E2").Value = "B5"
